{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Remove all text content, leaving a single empty paragraph behind (Word\n// always needs at least one paragraph in the body). Clear every paragraph's\n// text, then delete the extra (now-empty) paragraphs so only the first one\n// remains - mirroring the final state of a fully emptied body.\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  if (i === 0) {\n    paragraphs.items[i].insertText(\"\", \"Replace\");\n  } else {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The whole body held a single paragraph of instructional text (the \"Cuadro\n# de tiempos invertidos\" / \"Desviaci\u00f3n del tiempo...\" notes). The edit wipes\n# that content, leaving one empty paragraph before the section break - so\n# clear the entire story range in place.\n$d.Content.Delete()\n"}
